$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (465-491): date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44539, 1, 20, 323.4675723758693),
    @(44540, 7, 18, 291.1208151382824),
    @(44541, 0, 17, 274.9474365194889),
    @(44542, 2, 18, 291.1208151382824),
    @(44543, 3, 13, 210.2539220443151),
    @(44544, 4, 17, 274.9474365194889),
    @(44545, 0, 17, 274.9474365194889),
    @(44546, 1, 17, 274.9474365194889),
    @(44547, 2, 12, 194.0805434255216),
    @(44548, 1, 13, 210.2539220443151),
    @(44550, 5, 16, 258.7740579006955),
    @(44551, 4, 17, 274.9474365194889),
    @(44552, 0, 13, 210.2539220443151),
    @(44553, 0, 13, 210.2539220443151),
    @(44554, 4, 16, 258.7740579006955),
    @(44555, 6, 20, 323.4675723758693),
    @(44556, 0, 19, 307.2941937570758),
    @(44557, 7, 21, 339.6409509946628),
    @(44558, 5, 22, 355.8143296134562),
    @(44559, 8, 30, 485.2013585638039),
    @(44560, 8, 38, 614.5883875141517),
    @(44561, 11, 45, 727.802037845706),
    @(44562, 19, 58, 938.0559598900211),
    @(44563, 3, 61, 986.5760957464014),
    @(44564, 14, 68, 1099.789746077956),
    @(44565, 12, 75, 1213.00339640951),
    @(44566, 13, 80, 1293.870289503477)
)

$startRow = 465
$endRow = 491

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Copy the date-column style (s="2", YYYY-MM-DD HH:MM:SS format) from the last existing row down to the new rows
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0 | Out-Null
